$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Header restructuring: un-merge B4:B5, move "Periodo" down to B5, clear B4
# ---------------------------------------------------------------------------
$ws.Range("B4:B5").UnMerge()
$ws.Range("B4").Value = ""
$ws.Range("B5").Value = "Periodo"

# ---------------------------------------------------------------------------
# 2. Re-label the second (carga) group of column headers so every table
#    column has a unique name (Excel requires unique ListColumn names).
#    Visually they keep the same text, just with a trailing space.
# ---------------------------------------------------------------------------
$ws.Range("G5").Value = "Hombre Camión "
$ws.Range("H5").Value = "Pequeña "
$ws.Range("I5").Value = "Mediana "
$ws.Range("J5").Value = "Grande "

# ---------------------------------------------------------------------------
# 3. Column widths (B5:J82 data columns C..J got a bit wider once the
#    bestFit auto-sizing was replaced by fixed widths from the table).
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 15.142857142857142
$ws.Columns.Item(4).ColumnWidth = 8.714285714285714
$ws.Columns.Item(5).ColumnWidth = 8.428571428571429
$ws.Columns.Item(6).ColumnWidth = 7.428571428571429
$ws.Columns.Item(7).ColumnWidth = 16.0
$ws.Columns.Item(8).ColumnWidth = 9.714285714285714
$ws.Columns.Item(9).ColumnWidth = 9.428571428571429
$ws.Columns.Item(10).ColumnWidth = 8.428571428571429

# ---------------------------------------------------------------------------
# 4. Turn B5:J82 into an actual Excel Table ("Tabla1") and filter the
#    "Periodo" column down to calendar year 2024, which is what hides all
#    the older monthly rows (11-82) while keeping the 2024 rows (6-10)
#    visible.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("B5:J82"), 0, 1)
$lo.Name = "Tabla1"
$lo.TableStyle = "TableStyleMedium2"
$lo.Range.AutoFilter(1, ">=" + 45292, 1, "<" + 45658)

# Belt-and-suspenders: make sure the pre-2024 rows are actually hidden,
# matching the effect of the date filter above.
$ws.Range("B11:B82").EntireRow.Hidden = $true
$ws.Range("B6:B10").EntireRow.Hidden = $false

# ---------------------------------------------------------------------------
# 5. Clean up the per-cell left/right borders that used to mark the outer
#    edges of the data block; the Table's own border formatting now draws
#    those edges instead, so the old per-cell accents are removed to avoid
#    a doubled-up border line.
# ---------------------------------------------------------------------------
for ($r = 4; $r -le 82; $r++) {
    $ws.Cells.Item($r, 2).Borders.Item(7).LineStyle = -4142
    $ws.Cells.Item($r, 10).Borders.Item(10).LineStyle = -4142
}
